$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 5 data: date/time value in column A, zeros in B..M, "Random" label in N
$ws.Range("A5").Value = 42607.890289351853
$ws.Range("A5").NumberFormat = "m/d/yy h:mm"

$ws.Range("B5").Value = 17
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 0
$ws.Range("N5").Value = "Random"

$wb.Save()
